$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E2").Value = '[''Normal'']'
$ws.Range("D11").Value = '[1, 0, 1, 0, 1, 0, 0]'
$ws.Range("E11").Value = '[''Normal'', ''HardwareFault'', ''RegulationViolation'']'
$ws.Range("D24").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E24").Value = '[''Normal'']'
$ws.Range("D25").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E25").Value = '[''Normal'']'
$ws.Range("D26").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E26").Value = '[''SoftwareFault'']'
$ws.Range("D27").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E27").Value = '[''SoftwareFault'']'
$ws.Range("D29").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E29").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D35").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E35").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D38").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E38").Value = '[''Normal'']'
$ws.Range("D54").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E54").Value = '[''SoftwareFault'']'
$ws.Range("D56").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E56").Value = '[''Normal'']'
$ws.Range("D58").Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Range("E58").Value = '[''Normal'', ''ParamViolation'']'
$ws.Range("D61").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E61").Value = '[''Normal'']'
$ws.Range("D68").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E68").Value = '[''Normal'']'
$ws.Range("D71").Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Range("E71").Value = '[''Normal'', ''ParamViolation'']'
$ws.Range("D82").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E82").Value = '[''Normal'', ''SurroundingEnvironment'']'
$ws.Range("D83").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E83").Value = '[''Normal'']'
$ws.Range("D109").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E109").Value = '[''Normal'', ''SurroundingEnvironment'']'
$ws.Range("D113").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E113").Value = '[''Normal'', ''HardwareFault'']'
